$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the new
# header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-23
$values = @(
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(10, 10),
    @(5, 6),
    @(8, 8),
    @(10, 10),
    @(5, 5),
    @(6, 7),
    @(5, 7),
    @(6, 7),
    @(8, 8),
    @(4, 4),
    @(7, 7),
    @(5, 7),
    @(1, 3),
    @(1, 3),
    @(3, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
